$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 89

$ws.Cells.Item($row, 1).Value = "Kindergarden"
$ws.Cells.Item($row, 2).Value = "Kindergarden Den Haag Bezuidenhoutseweg"
$ws.Cells.Item($row, 3).Value = "KDV"

# Column D holds the report date as plain text (e.g. "2024-09-23"), not a
# real date serial. Force text formatting first so Excel doesn't
# autoconvert the string into a date value, then restore the default
# "Normal" style so no explicit style index gets attached to the cell.
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "2024-09-23"
$ws.Cells.Item($row, 4).Style = "Normal"

$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
